# Apply the update to "G L PURAM - TEACHERS TABS DATA.xlsx"
# - Fill in treasury-id (col L) / SD-card-no (col M) values for rows 15-17 and 19-27
# - Hide those same rows (13-15, 17-25 serial numbers) since the sheet uses AutoFilter
# - Hide columns J and K (no longer "best fit", just hidden helper columns)
# - Move the active selection from L20 to G65

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Newly filled-in data for column L (treasury id) / column M (SD card no.) ---
$rowData = @{
    15 = @{ L = "R9ZTA05GMKR"; M = 1019986616 }
    16 = @{ L = "R9ZT90G2BSZ"; M = 1507447064 }
    17 = @{ L = "R9ZT90LJPIA"; M = 1813702190 }
    19 = @{ L = "R9ZTB07226Y"; M = 1117231401 }
    20 = @{ L = "R9ZTA05JH4Y"; M = 1205958610 }
    21 = @{ L = "R9ZTB07Y8VE"; M = 1415131633 }
    22 = @{ L = "R9ZTB07Y3CB"; M = 1052159883 }
    23 = @{ L = "R9ZTA05GBSN"; M = 1301544669 }
    24 = @{ L = "R9ZTB093MXL"; M = 1736560875 }
    25 = @{ L = "R9ZTB07XTXV"; M = 1794724002 }
    26 = @{ L = "R9ZTB07XXNW"; M = 1719634398 }
    27 = @{ L = "R9ZTA05J05X"; M = 1917822828 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("L$r").Value2 = $vals.L
    $ws.Range("M$r").Value2 = $vals.M
}

# --- Hide the rows that now have their data filled in ---
foreach ($r in $rowData.Keys) {
    $ws.Rows.Item($r).Hidden = $true
}

# --- Hide helper columns J and K ---
$ws.Columns.Item(10).Hidden = $true
$ws.Columns.Item(11).Hidden = $true

# --- Update the active selection / view ---
[void]$ws.Range("G65").Select()

Write-Host "done"
